$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column C
$ws.Range("C1").Value = "Douban"

# New data for column C (rows 2-11)
$cValues = @(1, 245, 650, 1116, 1576, 2126, 2649, 3122, 3563, 3998)
for ($i = 0; $i -lt $cValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $cValues[$i]
}

# Update selection to match target (C12)
$ws.Range("C12").Select()
